$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25 (bus indices 0-23), columns B-F and I-N
# (column G stays 1, column H stays empty, per the original layout)
$data = @{
    2 = @{ "B"=1.02; "C"=1.040682962810709; "D"=1.044650579505955; "E"=1.054481022958382; "F"=1.062391396068169; "I"=1.044160574827285; "J"=1.0457682464687; "K"=1.047421032056998; "L"=1.057224115175593; "M"=1.065112870021561; "N"=1.01915482959033 }
    3 = @{ "B"=1.02; "C"=1.041556920506355; "D"=1.045325319583083; "E"=1.055347744267757; "F"=1.063344579854361; "I"=1.044411456751314; "J"=1.046288299394156; "K"=1.047907409356284; "L"=1.057903951514892; "M"=1.06588051273279; "N"=1.019329123768549 }
    4 = @{ "B"=1.02; "C"=1.042122916122461; "D"=1.045762237462051; "E"=1.05590942465319; "F"=1.063962349535457; "I"=1.044572718385198; "J"=1.046624637172694; "K"=1.048221762584779; "L"=1.058344061627675; "M"=1.066377597984379; "N"=1.019441801398388 }
    5 = @{ "B"=1.02; "C"=1.042360975840898; "D"=1.045945991837236; "E"=1.05614575791192; "F"=1.064222296785902; "I"=1.044640254537769; "J"=1.046765991508257; "K"=1.04835382809211; "L"=1.058529133027847; "M"=1.066586659600242; "N"=1.019489146157459 }
    6 = @{ "B"=1.02; "C"=1.042400953829992; "D"=1.045976849305916; "E"=1.056185451154301; "F"=1.064265956917858; "I"=1.04465157899136; "J"=1.046789723025254; "K"=1.048375997244672; "L"=1.058560210172417; "M"=1.066621767025586; "N"=1.019497094077746 }
    7 = @{ "B"=1.02; "C"=1.042126096637208; "D"=1.045764692507775; "E"=1.055912581754755; "F"=1.063965822034549; "I"=1.044573621822149; "J"=1.046626526122878; "K"=1.04822352760004; "L"=1.05834653436944; "M"=1.066380391135896; "N"=1.019442434119797 }
    8 = @{ "B"=1.02; "C"=1.040978219369041; "D"=1.044878544735816; "E"=1.054773757769902; "F"=1.062713321807941; "I"=1.044245583865913; "J"=1.045944035601471; "K"=1.047585480534606; "L"=1.057453824676187; "M"=1.065372221215492; "N"=1.019213754070303 }
    9 = @{ "B"=1.02; "C"=1.038959297665824; "D"=1.043319526598338; "E"=1.052773612133121; "F"=1.060513954497015; "I"=1.043659331204626; "J"=1.044740135965017; "K"=1.046458410569969; "L"=1.055882425892225; "M"=1.063598584953197; "N"=1.018810024849428 }
    10 = @{ "B"=1.02; "C"=1.037615972418149; "D"=1.042281944745957; "E"=1.051444717069088; "F"=1.059052978804896; "I"=1.043263019212566; "J"=1.043936749743083; "K"=1.045705247802732; "L"=1.054836026120661; "M"=1.062418184480019; "N"=1.018540381595308 }
    11 = @{ "B"=1.02; "C"=1.037034935990106; "D"=1.041833097784394; "E"=1.050870385973797; "F"=1.05842162950438; "I"=1.043090121376833; "J"=1.043588701362182; "K"=1.045378710924994; "L"=1.054383224004327; "M"=1.061907555175213; "N"=1.01842351187532 }
    12 = @{ "B"=1.02; "C"=1.036819209456296; "D"=1.041666442598998; "E"=1.050657219030233; "F"=1.058187309543609; "I"=1.043025705966549; "J"=1.04345939517601; "K"=1.045257359715717; "L"=1.054215078865937; "M"=1.061717959976923; "N"=1.018380084774702 }
    13 = @{ "B"=1.02; "C"=1.03686547915334; "D"=1.041702187665354; "E"=1.050702936562787; "F"=1.058237563289762; "I"=1.043039532034193; "J"=1.043487132929904; "K"=1.045283392699221; "L"=1.05425114447911; "M"=1.061758625362957; "N"=1.01838940077458 }
    14 = @{ "B"=1.02; "C"=1.037017101996336; "D"=1.041819320648093; "E"=1.05085276214879; "F"=1.058402256632255; "I"=1.043084800721471; "J"=1.043578013387042; "K"=1.045368681231686; "L"=1.054369324134508; "M"=1.061891881633795; "N"=1.018419922511538 }
    15 = @{ "B"=1.02; "C"=1.037110534628507; "D"=1.041891499010466; "E"=1.05094509657388; "F"=1.058503755029226; "I"=1.04311266663912; "J"=1.043634004485763; "K"=1.045421222277629; "L"=1.054442144618256; "M"=1.061973995233852; "N"=1.018438725788998 }
    16 = @{ "B"=1.02; "C"=1.037654547322194; "D"=1.042311742490668; "E"=1.05148285661377; "F"=1.059094906117548; "I"=1.04327446670614; "J"=1.043959844931386; "K"=1.045726910383412; "L"=1.054866083448672; "M"=1.062452083744887; "N"=1.01854813551996 }
    17 = @{ "B"=1.02; "C"=1.037995962287459; "D"=1.042575467086384; "E"=1.051820471994125; "F"=1.05946605876684; "I"=1.043375614147716; "J"=1.044164189427563; "K"=1.045918550805518; "L"=1.055132089194112; "M"=1.062752108689489; "N"=1.018616735544181 }
    18 = @{ "B"=1.02; "C"=1.038195164779478; "D"=1.042729334840008; "E"=1.052017502324921; "F"=1.059682667566606; "I"=1.043434487022331; "J"=1.044283362981581; "K"=1.046030291554909; "L"=1.05528727430851; "M"=1.06292715553252; "N"=1.018656737904405 }
    19 = @{ "B"=1.02; "C"=1.038263097975598; "D"=1.042781806787073; "E"=1.052084702337175; "F"=1.05975654617135; "I"=1.043454539976627; "J"=1.044323995145785; "K"=1.046068385467635; "L"=1.055340193211332; "M"=1.062986849975599; "N"=1.01867037581885 }
    20 = @{ "B"=1.02; "C"=1.03795932538244; "D"=1.042547167613082; "E"=1.051784238188003; "F"=1.059426224996855; "I"=1.043364774879671; "J"=1.044142266961627; "K"=1.045897993717604; "L"=1.05510354634175; "M"=1.062719913942924; "N"=1.01860937653255 }
    21 = @{ "B"=1.02; "C"=1.036972450180494; "D"=1.041784826054418; "E"=1.050808637675293; "F"=1.058353753254979; "I"=1.043071475556281; "J"=1.04355125205362; "K"=1.045343567556301; "L"=1.054334521898965; "M"=1.061852638894842; "N"=1.018410935072492 }
    22 = @{ "B"=1.02; "C"=1.036352520172802; "D"=1.041305897519097; "E"=1.050196195461268; "F"=1.057680555165549; "I"=1.042885947689131; "J"=1.043179510388857; "K"=1.044994626166003; "L"=1.05385127091305; "M"=1.061307784763642; "N"=1.018286071885699 }
    23 = @{ "B"=1.02; "C"=1.036681103518437; "D"=1.041559749510495; "E"=1.050520771478368; "F"=1.058037324715184; "I"=1.042984405359633; "J"=1.043376591275241; "K"=1.045179639568435; "L"=1.05410742583988; "M"=1.061596580454223; "N"=1.018352273108601 }
    24 = @{ "B"=1.02; "C"=1.037975879834697; "D"=1.042559954798392; "E"=1.051800610361092; "F"=1.059444223786902; "I"=1.043369673063751; "J"=1.044152172833393; "K"=1.045907282704062; "L"=1.05511644354031; "M"=1.062734461216872; "N"=1.018612701786783 }
    25 = @{ "B"=1.02; "C"=1.039480780795649; "D"=1.043722265983096; "E"=1.05328990512589; "F"=1.061081621782915; "I"=1.043811859963243; "J"=1.045051516148775; "K"=1.046750104222225; "L"=1.056288464363962; "M"=1.064056761760812; "N"=1.01891448676089 }
}

foreach ($row in $data.Keys) {
    $rowVals = $data[$row]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$row").Value = $rowVals[$col]
    }
}
